# Automatische test-sync: 2025-08-03 18:21:50
# Append a new test-mail row to the "Logs" sheet, append the matching
# category tally row to the "Dashboard" sheet, and extend the chart's
# series references + conditional formatting ranges to cover the new row.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- 1. Append new row 33 to "Logs" -----------------------------------
$logs.Cells.Item(33, 1).Value = "Wil je deze klant bellen?"
$logs.Cells.Item(33, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(33, 3).Value = "Testmail #5: Wil je deze klant bellen?"
$logs.Cells.Item(33, 4).Value = "Klantenservice / Contact"
$logs.Cells.Item(33, 5).Value = "Bedankt, we hebben dit doorgestuurd naar klantenservice@bedrijf.nl."
$logs.Cells.Item(33, 6).Value = "2025-08-03 18:21:26"
$logs.Cells.Item(33, 7).Value = "Ja"
$logs.Cells.Item(33, 8).Value = "Ja"
$logs.Cells.Item(33, 9).Value = "Nee"
$logs.Cells.Item(33, 10).Value = "Nee"

# --- 2. Append new row 8 to "Dashboard" --------------------------------
$dash.Cells.Item(8, 1).Value = "Klantenservice / Contact"
$dash.Cells.Item(8, 2).Value = 1

# --- 3. Extend conditional formatting ranges on "Logs" to row 33 -------
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "32")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "33")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count(); $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- 4. Extend the Dashboard chart series to include the new row -------
$co = $dash.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES('Dashboard'!`$B`$1,'Dashboard'!`$A`$2:`$A`$8,'Dashboard'!`$B`$2:`$B`$8,1)"
